# Apply updated "dSF" (column F) values for the rows whose data was repulled.
# Mapping of worksheet row number -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    5  = 0
    6  = 0
    9  = 0
    10 = -3
    13 = -5
    14 = 7
    16 = -1
    21 = -2
    29 = -7
    36 = -3
    41 = -2
    44 = -2
    45 = 7
    46 = -2
    50 = 3
    52 = 3
    58 = 2
    59 = 4
    71 = -1
    73 = -2
    76 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
